# Add five new custom styles (3 paragraph, 2 character) to the template's
# style sheet, matching the "elegant" CV template's contact/skill styling.
#
# wdStyleTypeParagraph = 1, wdStyleTypeCharacter = 2
# wdLineSpaceMultiple  = 5 (LineSpacing is then expressed in points per
#                            single line, e.g. 18 -> 1.5 lines -> w:line=360)
# Font.Color takes a BGR-packed integer (same as VBA's RGB()/wdColor),
# i.e. B*65536 + G*256 + R, for a target RRGGBB hex color.

$d = $word.ActiveDocument

# --- Contact Info (paragraph) ---------------------------------------------
$contactInfo = $d.Styles.Add("Contact Info", 1)
$contactInfo.ParagraphFormat.SpaceBefore = 0
$contactInfo.ParagraphFormat.SpaceAfter = 5.1
$contactInfo.ParagraphFormat.LineSpacingRule = 5
$contactInfo.ParagraphFormat.LineSpacing = 18
$contactInfo.Font.Name = "Liberation Serif"
$contactInfo.Font.Bold = $false
$contactInfo.Font.Color = 1191292
$contactInfo.Font.Size = 9

# --- Skill Category (paragraph) -------------------------------------------
$skillCategory = $d.Styles.Add("Skill Category", 1)
$skillCategory.ParagraphFormat.SpaceBefore = 0
$skillCategory.ParagraphFormat.SpaceAfter = 0
$skillCategory.ParagraphFormat.LineSpacingRule = 5
$skillCategory.ParagraphFormat.LineSpacing = 14.4
$skillCategory.Font.Name = "Liberation Serif"
$skillCategory.Font.Bold = $true
$skillCategory.Font.Color = 423897
$skillCategory.Font.Size = 10

# --- Skill Items (paragraph) ----------------------------------------------
$skillItems = $d.Styles.Add("Skill Items", 1)
$skillItems.ParagraphFormat.SpaceBefore = 0
$skillItems.ParagraphFormat.SpaceAfter = 0
$skillItems.ParagraphFormat.LineSpacingRule = 5
$skillItems.ParagraphFormat.LineSpacing = 14.4
$skillItems.Font.Name = "Liberation Serif"
$skillItems.Font.Bold = $false
$skillItems.Font.Color = 1191292
$skillItems.Font.Size = 10

# --- Skill Highlight (character) -------------------------------------------
$skillHighlight = $d.Styles.Add("Skill Highlight", 2)
$skillHighlight.Font.Name = "Liberation Serif"
$skillHighlight.Font.Bold = $true
$skillHighlight.Font.Color = 423897
$skillHighlight.Font.Size = 10

# --- Skill Level (character) ------------------------------------------------
$skillLevel = $d.Styles.Add("Skill Level", 2)
$skillLevel.Font.Name = "Liberation Serif"
$skillLevel.Font.Bold = $false
$skillLevel.Font.Color = 423897
$skillLevel.Font.Size = 10
